$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before row 98, pushing existing rows 98-116 down to 100-118.
$ws.Rows.Item(98).Resize(2).Insert()

# Fill in the two new rows (98 and 99) with values, copying the common template
# fields from the (now shifted) row 100 which retains the same shared values.

$rows = @(98, 99)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = 6
    $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100101
    $ws.Cells.Item($r, 8).Value = "Berries"
    $ws.Cells.Item($r, 9).Value = 100101008
    $ws.Cells.Item($r, 10).Value = "Mora"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = "Primera"
    $ws.Cells.Item($r, 17).Value = "$/bandeja 2 kilos"
    $ws.Cells.Item($r, 20).Value = 2
}

# Row 98 specific values
$ws.Cells.Item(98, 4).Value = 44943
$ws.Cells.Item(98, 13).Value = 150
$ws.Cells.Item(98, 14).Value = 4000
$ws.Cells.Item(98, 15).Value = 4000
$ws.Cells.Item(98, 16).Value = 4000
$ws.Cells.Item(98, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(98, 19).Value = 2000

# Row 99 specific values
$ws.Cells.Item(99, 4).Value = 44943
$ws.Cells.Item(99, 13).Value = 200
$ws.Cells.Item(99, 14).Value = 4000
$ws.Cells.Item(99, 15).Value = 4000
$ws.Cells.Item(99, 16).Value = 4000
$ws.Cells.Item(99, 18).Value = "Región del Maule"
$ws.Cells.Item(99, 19).Value = 2000

# Apply the date number format (same one used by column D elsewhere) to the new D cells
$ws.Cells.Item(98, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(99, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
